$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Cells.Item(18, 8).Value = 3459.327
$ws.Cells.Item(18, 9).Value = 3409.5098
$ws.Cells.Item(18, 11).Value = 3409.5098
$ws.Cells.Item(18, 13).Value = -3125.5098
# Row 76
$ws.Cells.Item(76, 8).Value = 4141.3335
$ws.Cells.Item(76, 10).Value = 3950
$ws.Cells.Item(76, 12).Value = 3950
$ws.Cells.Item(76, 14).Value = -4580
# Row 79
$ws.Cells.Item(79, 8).Value = 4141.3335
$ws.Cells.Item(79, 10).Value = 3950
$ws.Cells.Item(79, 12).Value = 3950
$ws.Cells.Item(79, 14).Value = -6134
# Row 86
$ws.Cells.Item(86, 8).Value = 57963.375
$ws.Cells.Item(86, 9).Value = 5367.6665
$ws.Cells.Item(86, 10).Value = 89520.8
$ws.Cells.Item(86, 11).Value = 5367.6665
$ws.Cells.Item(86, 12).Value = 89520.8
$ws.Cells.Item(86, 13).Value = -4244.6665
$ws.Cells.Item(86, 14).Value = -91766.8
# Row 89
$ws.Cells.Item(89, 8).Value = 57963.375
$ws.Cells.Item(89, 9).Value = 5367.6665
$ws.Cells.Item(89, 10).Value = 89520.8
$ws.Cells.Item(89, 11).Value = 26838.3325
$ws.Cells.Item(89, 12).Value = 447604
$ws.Cells.Item(89, 13).Value = -21222.3325
$ws.Cells.Item(89, 14).Value = -458836
# Row 92
$ws.Cells.Item(92, 8).Value = 778.2308
$ws.Cells.Item(92, 9).Value = 759.75
$ws.Cells.Item(92, 11).Value = 759.75
$ws.Cells.Item(92, 13).Value = 488.25
# Row 132
$ws.Cells.Item(132, 8).Value = 2922.547
$ws.Cells.Item(132, 9).Value = 3158.7144
$ws.Cells.Item(132, 11).Value = 9476.143199999999
$ws.Cells.Item(132, 13).Value = -6946.143199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Cells.Item(10, 8).Value = 9500
$ws.Cells.Item(10, 10).Value = 9500
$ws.Cells.Item(10, 12).Value = 9500
$ws.Cells.Item(10, 14).Value = -9780
# Row 80
$ws.Cells.Item(80, 8).Value = 1092
$ws.Cells.Item(80, 9).Value = 1101.5
$ws.Cells.Item(80, 10).Value = 1082.5
$ws.Cells.Item(80, 11).Value = 1101.5
$ws.Cells.Item(80, 12).Value = 1082.5
$ws.Cells.Item(80, 13).Value = -103.5
$ws.Cells.Item(80, 14).Value = -3078.5
# Row 83
$ws.Cells.Item(83, 8).Value = 1092
$ws.Cells.Item(83, 9).Value = 1101.5
$ws.Cells.Item(83, 10).Value = 1082.5
$ws.Cells.Item(83, 11).Value = 5507.5
$ws.Cells.Item(83, 12).Value = 5412.5
$ws.Cells.Item(83, 13).Value = -515.5
$ws.Cells.Item(83, 14).Value = -15396.5
# Row 99
$ws.Cells.Item(99, 8).Value = 21413.111
$ws.Cells.Item(99, 9).Value = 23121
$ws.Cells.Item(99, 10).Value = 7750
$ws.Cells.Item(99, 11).Value = 23121
$ws.Cells.Item(99, 12).Value = 7750
$ws.Cells.Item(99, 13).Value = -21623
$ws.Cells.Item(99, 14).Value = -10746
# Row 105
$ws.Cells.Item(105, 8).Value = 2470.1428
$ws.Cells.Item(105, 9).Value = 2072.5454
$ws.Cells.Item(105, 11).Value = 2072.5454
$ws.Cells.Item(105, 13).Value = -325.5454

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2199.1875
$ws.Cells.Item(16, 9).Value = 2417.5454
$ws.Cells.Item(16, 11).Value = 2417.5454
$ws.Cells.Item(16, 13).Value = -2130.5454
# Row 31
$ws.Cells.Item(31, 8).Value = 5693
$ws.Cells.Item(31, 9).Value = 4662.3335
$ws.Cells.Item(31, 10).Value = 6536.273
$ws.Cells.Item(31, 11).Value = 4662.3335
$ws.Cells.Item(31, 12).Value = 6536.273
$ws.Cells.Item(31, 13).Value = -4367.3335
$ws.Cells.Item(31, 14).Value = -7126.273
# Row 34
$ws.Cells.Item(34, 8).Value = 5693
$ws.Cells.Item(34, 9).Value = 4662.3335
$ws.Cells.Item(34, 10).Value = 6536.273
$ws.Cells.Item(34, 11).Value = 4662.3335
$ws.Cells.Item(34, 12).Value = 6536.273
$ws.Cells.Item(34, 13).Value = -4460.3335
$ws.Cells.Item(34, 14).Value = -6940.273
# Row 113
$ws.Cells.Item(113, 8).Value = 2199.1875
$ws.Cells.Item(113, 9).Value = 2417.5454
$ws.Cells.Item(113, 11).Value = 2417.5454
$ws.Cells.Item(113, 13).Value = -247.5454
# Row 134
$ws.Cells.Item(134, 8).Value = 3686164.2
$ws.Cells.Item(134, 10).Value = 3866.8333
$ws.Cells.Item(134, 12).Value = 11600.4999
$ws.Cells.Item(134, 14).Value = -16670.4999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 215.26315
$ws.Cells.Item(2, 9).Value = 212.33333
$ws.Cells.Item(2, 11).Value = 1273.99998
$ws.Cells.Item(2, 13).Value = -1160.99998
# Row 17
$ws.Cells.Item(17, 8).Value = 1073.5555
$ws.Cells.Item(17, 9).Value = 420.33334
$ws.Cells.Item(17, 10).Value = 2380
$ws.Cells.Item(17, 11).Value = 1261.00002
$ws.Cells.Item(17, 12).Value = 7140
$ws.Cells.Item(17, 13).Value = -1092.00002
$ws.Cells.Item(17, 14).Value = -7478
# Row 38
$ws.Cells.Item(38, 8).Value = 1205.125
$ws.Cells.Item(38, 9).Value = 259.5625
$ws.Cells.Item(38, 10).Value = 2150.6875
$ws.Cells.Item(38, 11).Value = 778.6875
$ws.Cells.Item(38, 12).Value = 6452.0625
$ws.Cells.Item(38, 13).Value = -431.6875
$ws.Cells.Item(38, 14).Value = -7146.0625
# Row 131
$ws.Cells.Item(131, 8).Value = 11629653
$ws.Cells.Item(131, 9).Value = 166667890
$ws.Cells.Item(131, 10).Value = 1785.3875
$ws.Cells.Item(131, 11).Value = 500003670
$ws.Cells.Item(131, 12).Value = 5356.1625
$ws.Cells.Item(131, 13).Value = -499998630
$ws.Cells.Item(131, 14).Value = -15436.1625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 322.83334
$ws.Cells.Item(2, 9).Value = 287.4
$ws.Cells.Item(2, 11).Value = 287.4
$ws.Cells.Item(2, 13).Value = -174.4
# Row 70
$ws.Cells.Item(70, 8).Value = 10499.333
$ws.Cells.Item(70, 9).Value = 5500
$ws.Cells.Item(70, 11).Value = 5500
$ws.Cells.Item(70, 13).Value = -5230
# Row 73
$ws.Cells.Item(73, 8).Value = 10499.333
$ws.Cells.Item(73, 9).Value = 5500
$ws.Cells.Item(73, 11).Value = 5500
$ws.Cells.Item(73, 13).Value = -4564
# Row 80
$ws.Cells.Item(80, 8).Value = 3061.4443
$ws.Cells.Item(80, 9).Value = 2592.8333
$ws.Cells.Item(80, 10).Value = 3998.6667
$ws.Cells.Item(80, 11).Value = 2592.8333
$ws.Cells.Item(80, 12).Value = 3998.6667
$ws.Cells.Item(80, 13).Value = -1594.8333
$ws.Cells.Item(80, 14).Value = -5994.6667
# Row 83
$ws.Cells.Item(83, 8).Value = 3061.4443
$ws.Cells.Item(83, 9).Value = 2592.8333
$ws.Cells.Item(83, 10).Value = 3998.6667
$ws.Cells.Item(83, 11).Value = 12964.1665
$ws.Cells.Item(83, 12).Value = 19993.3335
$ws.Cells.Item(83, 13).Value = -7972.166499999999
$ws.Cells.Item(83, 14).Value = -29977.3335
# Row 113
$ws.Cells.Item(113, 8).Value = 3060.6924
$ws.Cells.Item(113, 9).Value = 3042.1428
$ws.Cells.Item(113, 11).Value = 3042.1428
$ws.Cells.Item(113, 13).Value = -872.1428000000001
# Row 125
$ws.Cells.Item(125, 8).Value = 42936.25
$ws.Cells.Item(125, 10).Value = 42936.25
$ws.Cells.Item(125, 12).Value = 42936.25
$ws.Cells.Item(125, 14).Value = -47856.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 318.57144
$ws.Cells.Item(16, 9).Value = 318.57144
$ws.Cells.Item(16, 11).Value = 318.57144
$ws.Cells.Item(16, 13).Value = -148.57144
# Row 22
$ws.Cells.Item(22, 8).Value = 3082.6667
$ws.Cells.Item(22, 10).Value = 3446.6667
$ws.Cells.Item(22, 12).Value = 3446.6667
$ws.Cells.Item(22, 14).Value = -4036.6667
# Row 27
$ws.Cells.Item(27, 8).Value = 3082.6667
$ws.Cells.Item(27, 10).Value = 3446.6667
$ws.Cells.Item(27, 12).Value = 3446.6667
$ws.Cells.Item(27, 14).Value = -3660.6667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 24346.191
$ws.Cells.Item(122, 9).Value = 2828.2778
$ws.Cells.Item(122, 10).Value = 72761.5
$ws.Cells.Item(122, 11).Value = 8484.8334
$ws.Cells.Item(122, 12).Value = 218284.5
$ws.Cells.Item(122, 13).Value = -6034.8334
$ws.Cells.Item(122, 14).Value = -223184.5
